$d = $word.ActiveDocument

# --- Hunk 1 ------------------------------------------------------------
# Expand the abbreviated title " «Компьютерный клуб" (inside the "На тему:"
# paragraph) into the full title
# " «Разработка информационное подсистемы «Компьютерный клуб»" and drop the
# _GoBack bookmark that used to sit right after that run (between it and
# the following run holding the closing guillemet).
$d.Content.Find.Execute(
    "«Компьютерный клуб", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "«Разработка информационное подсистемы «Компьютерный клуб»", 2
) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Hunk 2 --------------------------------------------------------------
# Further down there is a run of four otherwise-identical empty paragraphs
# (pPr holding only <w:szCs w:val="28"/>) immediately before the
# "Выполнил студент группы П-46-21" paragraph. The first two of those four
# get merged into a single paragraph, and the _GoBack bookmark re-appears
# inside that merged (now-first) paragraph.
$count = $d.Paragraphs.Count
$anchor = 0
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Выполнил студент*") {
        $anchor = $i
        break
    }
}

$target = $anchor - 4

# Merge paragraph $target with the following one by deleting the
# paragraph mark that ends it.
$p1 = $d.Paragraphs($target)
$d.Range($p1.Range.End - 1, $p1.Range.End).Delete()

# Re-seat the bookmark inside the freshly merged paragraph.
$merged = $d.Paragraphs($target)
$d.Bookmarks.Add("_GoBack", $merged.Range) | Out-Null
